$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace the content of a Range with a literal OOXML fragment via
# Range.InsertXML, wrapped in the pkg:package envelope the host expects.
# ---------------------------------------------------------------------------
function Set-RangeXml($range, [string]$bodyXml) {
    $frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:compression="store">' + `
        '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + $bodyXml + '</w:body>' + `
        '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    [void]$range.InsertXML($frag)
}

# ---------------------------------------------------------------------------
# "Wzmocnienie +N / +1K" -> "Wzmocnienie +N / +1 UK" in the four price-list
# cells, split across separate runs ("Wzmocnienie +N / +1", " ", "U", "K").
# The 4th occurrence additionally carries the (relocated) "_GoBack" bookmark
# between the space-run and the "U" run.
# ---------------------------------------------------------------------------
$pPr = 'w:rsidR="0058252D" w:rsidRDefault="0058252D" w:rsidP="00116DC9"'

for ($n = 1; $n -le 3; $n++) {
    $old = "Wzmocnienie +$n / +1K"
    $r = $d.Content
    $ok = $r.Find.Execute($old)
    if (-not $ok) { throw "could not find '$old'" }

    $body = "<w:p $pPr>" + `
        "<w:r><w:t>Wzmocnienie +$n / +1</w:t></w:r>" + `
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
        '<w:r><w:t>U</w:t></w:r>' + `
        '<w:r><w:t>K</w:t></w:r>' + `
        '</w:p>'
    Set-RangeXml $r $body
}

# 4th cell - also relocates the "_GoBack" bookmark here, between " " and "U".
$old4 = "Wzmocnienie +4 / +1K"
$r4 = $d.Content
$ok4 = $r4.Find.Execute($old4)
if (-not $ok4) { throw "could not find '$old4'" }

$body4 = "<w:p $pPr>" + `
    '<w:r><w:t>Wzmocnienie +4 / +1</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
    '<w:r><w:t>U</w:t></w:r>' + `
    '<w:r><w:t>K</w:t></w:r>' + `
    '</w:p>'
Set-RangeXml $r4 $body4

# ---------------------------------------------------------------------------
# "tru" + bookmarkStart/End "_GoBack" + "cizny" -> single run "trucizny"
# (the old bookmark location is removed - it moved to the 4th cell above).
# ---------------------------------------------------------------------------
$rStart = $d.Content
$okStart = $rStart.Find.Execute("Eliksiry")
if (-not $okStart) { throw "could not find 'Eliksiry'" }
$start = $rStart.Start
$rEnd = $d.Content
$okEnd = $rEnd.Find.Execute("i bomby")
if (-not $okEnd) { throw "could not find 'i bomby'" }
$end = $rEnd.End

$paraTarget = $d.Range($start, $end)
$paraBody = '<w:p w:rsidR="004C3EB7" w:rsidRDefault="004C3EB7" w:rsidP="00483193">' + `
    '<w:pPr><w:pStyle w:val="Nagwek1"/></w:pPr>' + `
    '<w:r><w:t>Eliksiry</w:t></w:r>' + `
    '<w:r w:rsidR="006B4ED8"><w:t xml:space="preserve">, </w:t></w:r>' + `
    '<w:r w:rsidR="007B4E6A"><w:t>trucizny</w:t></w:r>' + `
    '<w:r w:rsidR="001E4A3A"><w:t>,</w:t></w:r>' + `
    '<w:r w:rsidR="006B4ED8"><w:t xml:space="preserve"> olejki</w:t></w:r>' + `
    '<w:r w:rsidR="001E4A3A"><w:t xml:space="preserve"> i bomby</w:t></w:r>' + `
    '</w:p>'
Set-RangeXml $paraTarget $paraBody
